# Apply the "page_metadata" edit:
#  - Rename header labels (F1/G1/H1)
#  - Normalize "languages" column values to proper-cased language names
#  - Fix a couple of typos / shorten a couple of long descriptions
#  - Replace two "website name" cells with their actual deployed URLs
#  - Re-fit row 7 height and column G width after the content changes
#  - Move the active selection to H16

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row relabeling ---------------------------------------------
$ws.Range("F1").Value = "repo"
$ws.Range("G1").Value = "languages"
$ws.Range("H1").Value = "libraries_tools"

# --- Column G ("languages") capitalization / wording fixes -------------
$ws.Range("G4").Value  = "Python"
$ws.Range("G5").Value  = "Python"
$ws.Range("G6").Value  = "Python"
$ws.Range("G7").Value  = "Python"
$ws.Range("G8").Value  = "Python, HTML, CSS"
$ws.Range("G9").Value  = "Python"
$ws.Range("G10").Value = "Python"
$ws.Range("G11").Value = "Python"
$ws.Range("G12").Value = "Python"
$ws.Range("G13").Value = "Python"
$ws.Range("G14").Value = "Python"
$ws.Range("G15").Value = "Python, CSS"
$ws.Range("G16").Value = "Python, JavaScript, HTML, CSS"
$ws.Range("G18").Value = "Python"
$ws.Range("G20").Value = "JavaScript, HTML, CSS"
$ws.Range("G21").Value = "JavaScript, HTML, CSS"
$ws.Range("G22").Value = "JavaScript, HTML, CSS"
$ws.Range("G23").Value = "Markdown"

# --- Misc content fixes in other columns --------------------------------
$ws.Range("I5").Value = "webs scraping, OOP, NLP"
$ws.Range("I6").Value = "PCA, k-means clustering"
$ws.Range("I7").Value = "resampling, ensemble methods"

# --- Replace displayed project-name text with the live site URLs -------
$ws.Range("J21").Value = "https://cdpeters.github.io/biodiversity-dashboard-plotly/"
$ws.Range("J22").Value = "https://cdpeters.github.io/dynamic-UFO-website-javascript/"

# --- Re-fit sizes now that the text in row 7 / column G changed --------
$ws.Rows.Item(7).AutoFit()
$ws.Columns.Item(7).ColumnWidth = 30.6

# --- Move active selection ----------------------------------------------
$ws.Range("H16").Select()
